$wb = $excel.ActiveWorkbook

# --- Sheet ALC: 18 cell update(s) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 160
$ws.Range("I2").Value = 161.66667
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 161.66667
$ws.Range("L2").Value = 150
$ws.Range("M2").Value = -48.66667000000001
$ws.Range("N2").Value = -376
$ws.Range("H100").Value = 4827.75
$ws.Range("I100").Value = 2414.8462
$ws.Range("J100").Value = 6918.933
$ws.Range("K100").Value = 2414.8462
$ws.Range("L100").Value = 6918.933
$ws.Range("M100").Value = -1873.8462
$ws.Range("N100").Value = -8000.933
$ws.Range("H141").Value = 1592.0476
$ws.Range("I141").Value = 1612.579
$ws.Range("K141").Value = 4837.737
$ws.Range("M141").Value = 342.2629999999999

# --- Sheet ARM: 42 cell update(s) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 27741.666
$ws.Range("J24").Value = 27741.666
$ws.Range("L24").Value = 27741.666
$ws.Range("N24").Value = -28489.666
$ws.Range("H61").Value = 882475.7
$ws.Range("I61").Value = 985472.9399999999
$ws.Range("K61").Value = 985472.9399999999
$ws.Range("M61").Value = -985260.9399999999
$ws.Range("H74").Value = 4447
$ws.Range("I74").Value = 4079.4055
$ws.Range("J74").Value = 5493.231
$ws.Range("K74").Value = 4079.4055
$ws.Range("L74").Value = 5493.231
$ws.Range("M74").Value = -3205.4055
$ws.Range("N74").Value = -7241.231
$ws.Range("H77").Value = 4447
$ws.Range("I77").Value = 4079.4055
$ws.Range("J77").Value = 5493.231
$ws.Range("K77").Value = 20397.0275
$ws.Range("L77").Value = 27466.155
$ws.Range("M77").Value = -16029.0275
$ws.Range("N77").Value = -36202.155
$ws.Range("H96").Value = 54672
$ws.Range("J96").Value = 54672
$ws.Range("L96").Value = 54672
$ws.Range("N96").Value = -60164
$ws.Range("H100").Value = 27741.666
$ws.Range("J100").Value = 27741.666
$ws.Range("L100").Value = 27741.666
$ws.Range("N100").Value = -29905.666
$ws.Range("H136").Value = 882475.7
$ws.Range("I136").Value = 985472.9399999999
$ws.Range("K136").Value = 2956418.82
$ws.Range("M136").Value = -2953868.82
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280

# --- Sheet BSM: 20 cell update(s) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 2000
$ws.Range("J16").Value = 2000
$ws.Range("L16").Value = 2000
$ws.Range("N16").Value = -2340
$ws.Range("H26").Value = 23198.2
$ws.Range("I26").Value = 23198.2
$ws.Range("K26").Value = 23198.2
$ws.Range("M26").Value = -22906.2
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H99").Value = 33547.824
$ws.Range("I99").Value = 35457.375
$ws.Range("K99").Value = 35457.375
$ws.Range("M99").Value = -33959.375

# --- Sheet CRP: 34 cell update(s) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 305416.47
$ws.Range("I31").Value = 621095.3
$ws.Range("J31").Value = 18435.727
$ws.Range("K31").Value = 621095.3
$ws.Range("L31").Value = 18435.727
$ws.Range("M31").Value = -620800.3
$ws.Range("N31").Value = -19025.727
$ws.Range("H34").Value = 305416.47
$ws.Range("I34").Value = 621095.3
$ws.Range("J34").Value = 18435.727
$ws.Range("K34").Value = 621095.3
$ws.Range("L34").Value = 18435.727
$ws.Range("M34").Value = -620893.3
$ws.Range("N34").Value = -18839.727
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H58").Value = 564010.25
$ws.Range("J58").Value = 3066.6
$ws.Range("L58").Value = 3066.6
$ws.Range("N58").Value = -3472.6
$ws.Range("H105").Value = 19706.4
$ws.Range("I105").Value = 23834.5
$ws.Range("K105").Value = 23834.5
$ws.Range("M105").Value = -22087.5
$ws.Range("H136").Value = 564010.25
$ws.Range("J136").Value = 3066.6
$ws.Range("L136").Value = 9199.799999999999
$ws.Range("N136").Value = -14299.8

# --- Sheet CUL: 36 cell update(s) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2994.5
$ws.Range("I64").Value = 1996
$ws.Range("K64").Value = 5988
$ws.Range("M64").Value = -5718
$ws.Range("H67").Value = 2994.5
$ws.Range("I67").Value = 1996
$ws.Range("K67").Value = 5988
$ws.Range("M67").Value = -5052
$ws.Range("H75").Value = 4626.8184
$ws.Range("I75").Value = 350
$ws.Range("J75").Value = 6230.625
$ws.Range("K75").Value = 1050
$ws.Range("L75").Value = 18691.875
$ws.Range("M75").Value = -52
$ws.Range("N75").Value = -20687.875
$ws.Range("H78").Value = 4626.8184
$ws.Range("I78").Value = 350
$ws.Range("J78").Value = 6230.625
$ws.Range("K78").Value = 3150
$ws.Range("L78").Value = 56075.625
$ws.Range("M78").Value = 1842
$ws.Range("N78").Value = -66059.625
$ws.Range("H107").Value = 1144.4242
$ws.Range("I107").Value = 244.33333
$ws.Range("J107").Value = 1481.9584
$ws.Range("K107").Value = 732.99999
$ws.Range("L107").Value = 4445.8752
$ws.Range("M107").Value = 1187.00001
$ws.Range("N107").Value = -8285.8752
$ws.Range("H114").Value = 7511
$ws.Range("I114").Value = 4222
$ws.Range("J114").Value = 8450.714
$ws.Range("K114").Value = 12666
$ws.Range("L114").Value = 25352.142
$ws.Range("M114").Value = -9412
$ws.Range("N114").Value = -31860.142

# --- Sheet GSM: 34 cell update(s) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 39900
$ws.Range("I5").Value = 39800
$ws.Range("J5").Value = 39950
$ws.Range("K5").Value = 39800
$ws.Range("L5").Value = 39950
$ws.Range("M5").Value = -39688
$ws.Range("N5").Value = -40174
$ws.Range("H9").Value = 9665.666999999999
$ws.Range("J9").Value = 16999.666
$ws.Range("L9").Value = 16999.666
$ws.Range("N9").Value = -17339.666
$ws.Range("H80").Value = 486671.4
$ws.Range("J80").Value = 35787
$ws.Range("L80").Value = 35787
$ws.Range("N80").Value = -37783
$ws.Range("H83").Value = 486671.4
$ws.Range("J83").Value = 35787
$ws.Range("L83").Value = 178935
$ws.Range("N83").Value = -188919
$ws.Range("H103").Value = 15000
$ws.Range("J103").Value = 15000
$ws.Range("L103").Value = 15000
$ws.Range("N103").Value = -17344
$ws.Range("H118").Value = 14054.5
$ws.Range("J118").Value = 14054.5
$ws.Range("L118").Value = 14054.5
$ws.Range("N118").Value = -17368.5
$ws.Range("H126").Value = 697195.5600000001
$ws.Range("I126").Value = 1113400.2
$ws.Range("J126").Value = 3521
$ws.Range("K126").Value = 3340200.6
$ws.Range("L126").Value = 10563
$ws.Range("M126").Value = -3337730.6
$ws.Range("N126").Value = -15503

# --- Sheet LTW: 19 cell update(s) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 36787.5
$ws.Range("J64").Value = 36787.5
$ws.Range("L64").Value = 36787.5
$ws.Range("N64").Value = -37237.5
$ws.Range("H67").Value = 36787.5
$ws.Range("J67").Value = 36787.5
$ws.Range("L67").Value = 36787.5
$ws.Range("N67").Value = -38347.5
$ws.Range("H93").Value = 1090.9375
$ws.Range("I93").Value = 1023.6667
$ws.Range("K93").Value = 1023.6667
$ws.Range("M93").Value = 224.3333
$ws.Range("H122").Value = 5468.375
$ws.Range("I122").Value = 5249.9287
$ws.Range("J122").Value = 6997.5
$ws.Range("K122").Value = 15749.7861
$ws.Range("L122").Value = 20992.5
$ws.Range("M122").Value = -13299.7861
$ws.Range("N122").Value = -25892.5

# --- Sheet WVR: 22 cell update(s) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H113").Value = 2458.1667
$ws.Range("I113").Value = 271
$ws.Range("J113").Value = 3551.75
$ws.Range("K113").Value = 813
$ws.Range("L113").Value = 10655.25
$ws.Range("M113").Value = 1357
$ws.Range("N113").Value = -14995.25
$ws.Range("H126").Value = 3207.6924
$ws.Range("I126").Value = 2733.8096
$ws.Range("J126").Value = 5198
$ws.Range("K126").Value = 8201.4288
$ws.Range("L126").Value = 15594
$ws.Range("M126").Value = -5731.4288
$ws.Range("N126").Value = -20534
$ws.Range("H136").Value = 15291.25
$ws.Range("I136").Value = 20265.125
$ws.Range("K136").Value = 60795.375
$ws.Range("M136").Value = -58245.375
